$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "  -3.21%  "
$ws.Range("E3").Value = "  -5.07%  "
$ws.Range("E4").Value = "  +0.23%  "
$ws.Range("E5").Value = "  -5.23%  "
$ws.Range("E6").Value = "  -6.82%  "
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("E8").Value = "  -5.27%  "
$ws.Range("E9").Value = "  -6.35%  "
$ws.Range("E10").Value = "  -7.12%  "
$ws.Range("E11").Value = "  -10.89%  "
$ws.Range("E12").Value = "  -9.91%  "
$ws.Range("E13").Value = "  -4.19%  "
$ws.Range("E14").Value = "  -1.14%  "
$ws.Range("E15").Value = "  -7.90%  "
$ws.Range("E16").Value = "  -4.21%  "
$ws.Range("E18").Value = "  -10.20%  "
$ws.Range("E19").Value = "  -7.20%  "
$ws.Range("E20").Value = "  -6.51%  "
$ws.Range("E21").Value = "  -9.36%  "
$ws.Range("E22").Value = "  -9.92%  "
$ws.Range("E23").Value = "  -0.22%  "
$ws.Range("E24").Value = "  -7.42%  "
$ws.Range("E25").Value = "  -6.52%  "
$ws.Range("E26").Value = "  -3.10%  "
$ws.Range("E27").Value = "  -0.35%  "
$ws.Range("E28").Value = "  -11.22%  "
$ws.Range("E30").Value = "  -8.23%  "
$ws.Range("E31").Value = "  -2.75%  "
$ws.Range("E32").Value = "  -7.91%  "
$ws.Range("E33").Value = "  -7.47%  "
$ws.Range("E34").Value = "  -6.49%  "
$ws.Range("E35").Value = "  -3.89%  "
$ws.Range("E36").Value = "  -8.01%  "
$ws.Range("E37").Value = "  -8.41%  "
$ws.Range("E38").Value = "  -10.48%  "
$ws.Range("E39").Value = "  -4.45%  "
$ws.Range("E40").Value = "  -7.50%  "
$ws.Range("E41").Value = "  -11.61%  "
$ws.Range("E42").Value = "  -11.26%  "
$ws.Range("E44").Value = "  -8.18%  "
$ws.Range("E45").Value = "  -5.72%  "
$ws.Range("E46").Value = "  -8.22%  "
$ws.Range("E47").Value = "  -3.15%  "
$ws.Range("E48").Value = "  -9.67%  "
$ws.Range("E49").Value = "  -5.89%  "
$ws.Range("E50").Value = "  -6.42%  "
$ws.Range("E51").Value = "  -7.42%  "
